$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header in column H, styled like the other header cells (G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "Label" column values (0 = Control rows, 1 = MDD rows)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1

# Refit values for the MDD 2 row (Control 33 changed between the two runs)
$ws.Range("D3").Value = 0.2576901534734823
$ws.Range("E3").Value = 0.2576901534734823
